# The deck originally has 7 slides:
#   1: Small Genome Assembly (title)
#   2: Small genome assembly (aim/method questions)
#   3: Small genome assembly (aim/method - what I would choose)
#   4: Reference guided assembly
#   5: Small genome assembly (reference guided assembly figure)
#   6: Small genome assembly (additional steps / circularization)
#   7: Small genome assembly (Illumina/Pacbio/Nanopore comparison)
#
# The commit removes slides 2-6, keeping only the original slide 1 and
# slide 7 (which becomes the new slide 2).

$p = $ppt.ActivePresentation

# Delete from the end backwards so indices of not-yet-deleted slides
# remain stable while iterating.
$p.Slides.Item(6).Delete()
$p.Slides.Item(5).Delete()
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
